$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Juice_test")

# --- Unhide columns D:E (item / preconditions columns) ---
$ws.Range("D1:E1").EntireColumn.Hidden = $false

# --- Fill in the previously-hidden Item/Preconditions cells on row 13 ---
$ws.Range("D12").Copy()
$ws.Range("D13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D13").Value2 = "Exercise2"

$ws.Range("E12").Copy()
$ws.Range("E13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E13").Value2 = "order exist"

# --- Add the new test case row (row 14) ---
$ws.Cells.Item(14, 1).Value2 = 13
$ws.Cells.Item(14, 2).Value2 = "XXX buys 3 cartons"
$ws.Cells.Item(14, 3).Value2 = "member must be boolean"
$ws.Cells.Item(14, 4).Value2 = "Exercise2"
$ws.Cells.Item(14, 5).Value2 = "order exist"
$ws.Cells.Item(14, 6).Value2 = 3
$ws.Cells.Item(14, 7).Value2 = """XXX"""
$ws.Cells.Item(14, 8).Value2 = "Value error"

# Match the Item/Preconditions styling used elsewhere in the table
$ws.Range("D12").Copy()
$ws.Range("D14").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E12").Copy()
$ws.Range("E14").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# --- Update the active selection to reflect where the user ended up ---
$ws.Range("C18").Select() | Out-Null
